# Append a new row of profit data (row 12) for 12/06/2025 to the worksheet,
# matching the run performed on 2025-12-06.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# Column A holds a date-like string ("12/06/2025") that must stay literal
# text (as the other rows' dates are stored), not get auto-converted into
# an Excel date serial number. Force text entry, then clear the formatting
# so no extra style gets attached to the cell (matching the plain,
# style-less cells used for every other data row).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "12/06/2025"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = 13382.72
$ws.Cells.Item($row, 3).Value = 0.1706518523656737
$ws.Cells.Item($row, 4).Value = 0.8293481476343263
$ws.Cells.Item($row, 5).Value = -88.65000000000001
$ws.Cells.Item($row, 6).Value = -20.38
$ws.Cells.Item($row, 7).Value = -19392.82
$ws.Cells.Item($row, 8).Value = -63.6
$ws.Cells.Item($row, 9).Value = -519.0599999999999
$ws.Cells.Item($row, 10).Value = -18.52
